$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The 3 archival-record rows ("aspect"/"blow"/"feather") that used to sit at
# the top of the data block (rows 17-19) are being dropped entirely, so
# every row below shifts up by three.
$ws.Rows("17:19").Delete()

# The 8 trailing metadata columns (old BC:BJ) are being dropped entirely, so
# the remaining BK:BN block (Entidades / N. area desc. fis. / N. sobre
# ilustrac. / Otras notas) shifts left into BC:BF.
$ws.Columns("BC:BJ").Delete()

# The bottommost formatted rows (old 35:37) that used to carry the sheet's
# row style no longer apply to real data; that formatting rolls down to the
# very end of the worksheet, leaving three empty, short (12.8pt) rows there.
$ws.Rows("1048574:1048576").RowHeight = 12.8

# Selection moved to A16.
$ws.Range("A16").Select()
